$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.394.40'
$ws.Range("E2").Value = '  -2.36%  '

# Row 3
$ws.Range("D3").Value = '1.845.00'
$ws.Range("E3").Value = '  -2.08%  '

# Row 4
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -8.32%  '

# Row 6
$ws.Range("E6").Value = '  +0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5253'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.34%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3224'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -9.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06756'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.83%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.91'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7700'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.30%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07718'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.41%  '

# Row 13
$ws.Range("D13").Value = '1.793.11'
$ws.Range("E13").Value = '  -4.86%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.15%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.028'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '

# Row 17
$ws.Range("E17").Value = '  -3.03%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007898'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.44%  '

# Row 20
$ws.Range("D20").Value = '26.441.39'
$ws.Range("E20").Value = '  -2.30%  '

# Row 21
$ws.Range("D21").Value = '2.120.48'
$ws.Range("E21").Value = '  -0.52%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.538'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.513'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -6.93%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.925'
$ws.Range("D24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = '  -2.60%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '143.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.14%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.649'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.30%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.59'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.14%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.193'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.18%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.154'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.64%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08793'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.07%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04804'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.53%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.132'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.15%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.853'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.68%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6900'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.99%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.111'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01791'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.217'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.22%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4916'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.00%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.85%  '

# Row 42
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8994'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.88%  '

# Row 43
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.194'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.12%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.14%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.785'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.36%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4197'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1265'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.83%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.104'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.05%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05884'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.13%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.40%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.98%  '
